# Refresh the cryptos list (price + 1h volume change columns) with the
# latest scraped values, mirroring the GitHub Actions update commit.
# Numeric-looking "Price" strings are forced to text (matching the
# inline-string storage already used in the sheet) by briefly flipping the
# cell to the "@" text format before assigning, then resetting the style so
# no visible formatting change is introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.012.62'
$ws.Range("E2").Value = '  +0.25%  '
$ws.Range("D3").Value = '2.415.87'
$ws.Range("E3").Value = '  -0.04%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '552.29'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.36%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '137.01'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.05%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("E8").Value = '  +3.82%  '
$ws.Range("E9").Value = '  -1.81%  '
$ws.Range("E10").Value = '  -2.73%  '
$ws.Range("E11").Value = '  -0.83%  '
$ws.Range("E12").Value = '  -1.60%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '25.25'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.02%  '
$ws.Range("D14").Value = '2.845.92'
$ws.Range("E14").Value = '  -0.07%  '
$ws.Range("D15").Value = '59.916.56'
$ws.Range("E15").Value = '  +0.28%  '
$ws.Range("E16").Value = '  -1.49%  '
$ws.Range("D17").Value = '2.386.61'
$ws.Range("E17").Value = '  -0.99%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.31'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.93%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.41'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.19%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '328.83'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.42%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.00'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.09%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '65.92'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.97%  '
$ws.Range("E24").Value = '  +2.83%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '8.62'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.76%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.05%  '
$ws.Range("E27").Value = '  +0.61%  '
$ws.Range("D28").Value = '0.0₃0773'
$ws.Range("E28").Value = '  -1.73%  '
$ws.Range("E29").Value = '  -1.99%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '169.30'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.91%  '
$ws.Range("E31").Value = '  -3.79%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '18.61'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.44%  '
$ws.Range("E33").Value = '  -1.29%  '
$ws.Range("E35").Value = '  -0.31%  '
$ws.Range("E36").Value = '  +0.08%  '
$ws.Range("E37").Value = '  -1.38%  '
$ws.Range("E38").Value = '  -1.94%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '322.15'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.13%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.404'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.49%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '140.10'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.07%  '
$ws.Range("E43").Value = '  +0.52%  '
$ws.Range("E44").Value = '  +1.92%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0514'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.83%  '
$ws.Range("E46").Value = '  +1.27%  '
$ws.Range("E47").Value = '  -1.36%  '
$ws.Range("E48").Value = '  -8.16%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '11.05'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.06%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.56'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.42%  '
$ws.Range("E51").Value = '  -0.92%  '
